# Update "想去人数" (F column) counts on three sheets to reflect refreshed
# scrape data (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 172
$ws1.Range("F4").Value  = 409
$ws1.Range("F7").Value  = 1085
$ws1.Range("F8").Value  = 364
$ws1.Range("F9").Value  = 187
$ws1.Range("F13").Value = 367
$ws1.Range("F14").Value = 777
$ws1.Range("F15").Value = 152
$ws1.Range("F16").Value = 715
$ws1.Range("F17").Value = 270
$ws1.Range("F21").Value = 255
$ws1.Range("F26").Value = 461

# --- Sheet "演出" -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 37

# --- Sheet "全部类型" ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 172
$ws4.Range("F6").Value  = 409
$ws4.Range("F9").Value  = 1085
$ws4.Range("F10").Value = 364
$ws4.Range("F11").Value = 187
$ws4.Range("F16").Value = 37
$ws4.Range("F20").Value = 367
$ws4.Range("F21").Value = 777
$ws4.Range("F22").Value = 152
$ws4.Range("F23").Value = 715
$ws4.Range("F24").Value = 270
$ws4.Range("F30").Value = 255
$ws4.Range("F38").Value = 461
